$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.151.16'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '2.271.55'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.73'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.52'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.16%  '
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('B8').Value = 'BinanceUSD'
$ws.Range('C8').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '50.05'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +4,902.00%  '
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.489'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.16%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '32.89'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.63%  '
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0805'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.113'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.74%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.69'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.623.27'
$ws.Range('E15').Value = '  +0.65%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.31'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.56%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.277.24'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('B18').Value = 'Polygon'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.785'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.61%  '
$ws.Range('D19').Value = '42.021.46'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.68'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.24%  '
$ws.Range('E21').Value = '  +1.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.98'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.20'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '244.20'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.30%  '
$ws.Range('E26').Value = '  +2.25%  '
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.04'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('E29').Value = '  +0.52%  '
$ws.Range('E31').Value = '  +3.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '160.12'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.33'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0743'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('E37').Value = '  +3.88%  '
$ws.Range('E38').Value = '  -1.13%  '
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.117'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.79'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.01'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.73'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').Value = '2.014.66'
$ws.Range('E44').Value = '  -2.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.25'
$ws.Range('D45').ClearFormats()
$ws.Range('E46').Value = '  +1.65%  '
$ws.Range('E47').Value = '  +2.09%  '
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.22'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.12%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.72'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.86%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.52'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.10%  '
